$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = -38915.87
$ws.Range("C2").Value = -38915.87
$ws.Range("B3").Value = 57327.05
$ws.Range("C3").Value = 57327.05
$ws.Range("B4").Value = -36353.05
$ws.Range("C4").Value = -36353.05
$ws.Range("B5").Value = 452802.19
$ws.Range("C5").Value = 452802.19
$ws.Range("B6").Value = -11866.61
$ws.Range("C6").Value = -11866.61
$ws.Range("B7").Value = -42160.42
$ws.Range("C7").Value = -42160.42
$ws.Range("B8").Value = -11891.18
$ws.Range("C8").Value = -11891.18
$ws.Range("B9").Value = -57936.9
$ws.Range("C9").Value = -57936.9
$ws.Range("B10").Value = -28564.63
$ws.Range("C10").Value = -28564.63
$ws.Range("B11").Value = 196582.08
$ws.Range("C11").Value = 196582.08
$ws.Range("B12").Value = -196582.08
$ws.Range("C12").Value = -196582.08
$ws.Range("B13").Value = -18135.67
$ws.Range("C13").Value = -18135.67
$ws.Range("B14").Value = 81716.10000000001
$ws.Range("C14").Value = 81716.10000000001
$ws.Range("B15").Value = 346021.01
$ws.Range("C15").Value = 346021.01

$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = -42885.77
$ws.Range("C2").Value = -42885.77
$ws.Range("B3").Value = 69089.91
$ws.Range("C3").Value = 69089.91
$ws.Range("B4").Value = -30457.38
$ws.Range("C4").Value = -30457.38
$ws.Range("B5").Value = 439948.04
$ws.Range("C5").Value = 439948.04
$ws.Range("B6").Value = -20344.57
$ws.Range("C6").Value = -20344.57
$ws.Range("B7").Value = -40308.66
$ws.Range("C7").Value = -40308.66
$ws.Range("B8").Value = -20277.56
$ws.Range("C8").Value = -20277.56
$ws.Range("B9").Value = -44474.84
$ws.Range("C9").Value = -44474.84
$ws.Range("B10").Value = -29994.23
$ws.Range("C10").Value = -29994.23
$ws.Range("B11").Value = 188143.44
$ws.Range("C11").Value = 188143.44
$ws.Range("B12").Value = -188143.44
$ws.Range("C12").Value = -188143.44
$ws.Range("B13").Value = -16482.58
$ws.Range("C13").Value = -16482.58
$ws.Range("B14").Value = 76255.21000000001
$ws.Range("C14").Value = 76255.21000000001
$ws.Range("B15").Value = 340067.57
$ws.Range("C15").Value = 340067.57

$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = -38293.06
$ws.Range("C2").Value = -38293.06
$ws.Range("B3").Value = 46505.97
$ws.Range("C3").Value = 46505.97
$ws.Range("B4").Value = -38819.41
$ws.Range("C4").Value = -38819.41
$ws.Range("B5").Value = 521104.54
$ws.Range("C5").Value = 521104.54
$ws.Range("B6").Value = -15438.74
$ws.Range("C6").Value = -15438.74
$ws.Range("B7").Value = -31167.05
$ws.Range("C7").Value = -31167.05
$ws.Range("B8").Value = -16132.08
$ws.Range("C8").Value = -16132.08
$ws.Range("B9").Value = -54028.1
$ws.Range("C9").Value = -54028.1
$ws.Range("B10").Value = -27644.43
$ws.Range("C10").Value = -27644.43
$ws.Range("B11").Value = 185434.74
$ws.Range("C11").Value = 185434.74
$ws.Range("B12").Value = -185434.74
$ws.Range("C12").Value = -185434.74
$ws.Range("B13").Value = -15629.56
$ws.Range("C13").Value = -15629.56
$ws.Range("B14").Value = 76187.89999999999
$ws.Range("C14").Value = 76187.89999999999
$ws.Range("B15").Value = 406645.98
$ws.Range("C15").Value = 406645.98

$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = -31954.07
$ws.Range("C2").Value = -31954.07
$ws.Range("B3").Value = 68042.99000000001
$ws.Range("C3").Value = 68042.99000000001
$ws.Range("B4").Value = -38946.09
$ws.Range("C4").Value = -38946.09
$ws.Range("B5").Value = 400785.61
$ws.Range("C5").Value = 400785.61
$ws.Range("B6").Value = -15008.15
$ws.Range("C6").Value = -15008.15
$ws.Range("B7").Value = -37792.51
$ws.Range("C7").Value = -37792.51
$ws.Range("B8").Value = -18620.19
$ws.Range("C8").Value = -18620.19
$ws.Range("B9").Value = -73444.08
$ws.Range("C9").Value = -73444.08
$ws.Range("B10").Value = -29930.08
$ws.Range("C10").Value = -29930.08
$ws.Range("B11").Value = 205682.74
$ws.Range("C11").Value = 205682.74
$ws.Range("B12").Value = -205682.74
$ws.Range("C12").Value = -205682.74
$ws.Range("B13").Value = -19907.67
$ws.Range("C13").Value = -19907.67
$ws.Range("B14").Value = 68580.99000000001
$ws.Range("C14").Value = 68580.99000000001
$ws.Range("B15").Value = 271806.75
$ws.Range("C15").Value = 271806.75

$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = -28973.95
$ws.Range("C2").Value = -28973.95
$ws.Range("B3").Value = 97647.47
$ws.Range("C3").Value = 97647.47
$ws.Range("B4").Value = -34430.07
$ws.Range("C4").Value = -34430.07
$ws.Range("B5").Value = 439761.29
$ws.Range("C5").Value = 439761.29
$ws.Range("B6").Value = -18061.25
$ws.Range("C6").Value = -18061.25
$ws.Range("B7").Value = -39080.4
$ws.Range("C7").Value = -39080.4
$ws.Range("B8").Value = -13311.51
$ws.Range("C8").Value = -13311.51
$ws.Range("B9").Value = -49862.93
$ws.Range("C9").Value = -49862.93
$ws.Range("B10").Value = -24965.85
$ws.Range("C10").Value = -24965.85
$ws.Range("B11").Value = 189920.85
$ws.Range("C11").Value = 189920.85
$ws.Range("B12").Value = -189920.85
$ws.Range("C12").Value = -189920.85
$ws.Range("B13").Value = -12846.96
$ws.Range("C13").Value = -12846.96
$ws.Range("B14").Value = 110584.34
$ws.Range("C14").Value = 110584.34
$ws.Range("B15").Value = 426460.18
$ws.Range("C15").Value = 426460.18
